$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 2-39: price / volume updates ---
$ws.Range("D2").Value = "65.552.81"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "3.713.25"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "'410.85"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "'133.47"
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("D7").Value = "3.702.94"
$ws.Range("E7").Value = "  +4.48%  "
$ws.Range("D8").Value = "'0.627"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "'0.737"
$ws.Range("E10").Value = "  -5.05%  "
$ws.Range("D11").Value = "'0.166"
$ws.Range("E11").Value = "  -6.59%  "
$ws.Range("D12").Value = "'0.0000340"
$ws.Range("E12").Value = "  +6.29%  "
$ws.Range("D13").Value = "'42.46"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").Value = "'10.02"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "4.311.65"
$ws.Range("E15").Value = "  +4.78%  "
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "3.723.06"
$ws.Range("E17").Value = "  +3.68%  "
$ws.Range("D18").Value = "'20.16"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").Value = "'12.98"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").Value = "'1.09"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").Value = "65.937.18"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").Value = "'423.42"
$ws.Range("E22").Value = "  -6.59%  "
$ws.Range("D23").Value = "'15.06"
$ws.Range("E23").Value = "  +14.35%  "
$ws.Range("D24").Value = "'87.12"
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("E25").Value = "  -5.48%  "
$ws.Range("D26").Value = "'36.42"
$ws.Range("E26").Value = "  +5.10%  "
$ws.Range("D27").Value = "'3.21"
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("E28").Value = "  -4.83%  "
$ws.Range("D29").Value = "'5.17"
$ws.Range("E29").Value = "  +7.04%  "
$ws.Range("D30").Value = "'12.59"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").Value = "'2.72"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").Value = "'6.99"
$ws.Range("E33").Value = "  -4.69%  "
$ws.Range("D34").Value = "'41.46"
$ws.Range("E34").Value = "  +6.22%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "'55.80"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'0.0472"
$ws.Range("E38").Value = "  -5.47%  "
$ws.Range("D39").Value = "'2.97"
$ws.Range("E39").Value = "  +25.56%  "

# --- rows 40-42: coin rotation (PEPE / Stellar / FirstDigitalUSD) ---
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.141"
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.996"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0663"
$ws.Range("E42").Value = "  -17.68%  "

# --- row 43 ---
$ws.Range("D43").Value = "'3.41"
$ws.Range("E43").Value = "  +4.70%  "

# --- rows 44-47: coin rotation (ApeXProtocol / EnergySwap / Monero / ARBITRUM) ---
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'27.24"
$ws.Range("E44").Value = "  +25.07%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'145.25"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'2.09"
$ws.Range("E46").Value = "  +5.08%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.08"
$ws.Range("E47").Value = "  +19.17%  "

# --- rows 48-51 ---
$ws.Range("E48").Value = "  -6.33%  "
$ws.Range("D49").Value = "'4.23"
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("E50").Value = "  -8.54%  "
$ws.Range("E51").Value = "  -5.57%  "
